# "leak before break satisfied" -- add a new "Leak Before Break" section
# below the existing matrix (rows 64-67): fracture-toughness input,
# critical-crack-size and FOS-on-critical-crack-size derived rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 64: section header (reuse the bold "section title" style used
#     elsewhere, e.g. A6 "Section Lengths") ---
$ws.Range("A64").Value = "Leak Before Break"
$ws.Range("A64").Font.Bold = $true

# --- Row 65: fracture toughness input row (plain literal numbers) ---
$ws.Range("A65").Value = "6061-T6 Fracture Toughness"
$ws.Range("B65").Value = "psi*sqrt(in)"
$ws.Range("C65:K65").Value = 26400

# --- Row 66: critical crack size, derived from fracture toughness and
#     the already-computed von Mises stress row (row 57) ---
$ws.Range("A66").Value = "Critical crack size"
$ws.Range("B66").Value = "in"
$ws.Range("C66").Formula = "=(C65/C57)^2/PI()"
$ws.Range("D66:K66").Formula = "=(D65/D57)^2/PI()"
$ws.Range("C66:K66").NumberFormat = "0.0000"

# --- Row 67: FOS on critical crack size vs thickness ---
$ws.Range("A67").Value = "FOS Cc vs thickness"
$ws.Range("B67").Value = "ratio"
$ws.Range("C67").Formula = "=C66/C3"
$ws.Range("D67:K67").Formula = "=D66/D3"
$ws.Range("C67:K67").NumberFormat = "0.0"

# --- View state: land the selection where the author left off ---
$ws.Range("L65").Select()
